$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2379653333333333
$ws.Range("H2").Value = 0.713896
$ws.Range("I2").Value = 0.0004000853538884766
$ws.Range("J2").Value = 0.0004000853538884766
$ws.Range("M2").Value = 1.675840666666667
$ws.Range("N2").Value = 5.027522
$ws.Range("O2").Value = 0.03808750486214892
$ws.Range("P2").Value = 0.03808750486214892
$ws.Range("Q2").Value = 0.3987919828568889
$ws.Range("R2").Value = 3.589127845712
$ws.Range("S2").Value = 0.00001523825286150192
$ws.Range("T2").Value = 0.00001523825286150192

# Row 3
$ws.Range("G3").Value = 0.2379653333333333
$ws.Range("H3").Value = 0.713896
$ws.Range("I3").Value = 0.0004000853538884766
$ws.Range("J3").Value = 0.0004000853538884766
$ws.Range("O3").Value = 0.096040539564286
$ws.Range("P3").Value = 0.09604053956428601
$ws.Range("Q3").Value = 1.005584307664889
$ws.Range("R3").Value = 9.050258768984
$ws.Range("S3").Value = 0.00003842441325921761
$ws.Range("T3").Value = 0.00003842441325921761

# Row 4
$ws.Range("G4").Value = 0.2379653333333333
$ws.Range("H4").Value = 0.713896
$ws.Range("I4").Value = 0.0004000853538884766
$ws.Range("J4").Value = 0.0004000853538884766
$ws.Range("M4").Value = 38.098149
$ws.Range("N4").Value = 114.294447
$ws.Range("O4").Value = 0.8658719555735651
$ws.Range("P4").Value = 0.865871955573565
$ws.Range("Q4").Value = 9.066038726167999
$ws.Range("R4").Value = 81.594348535512
$ws.Range("S4").Value = 0.0003464226877677571
$ws.Range("T4").Value = 0.0003464226877677571

# Row 5
$ws.Range("I5").Value = 0.9842542228653065
$ws.Range("J5").Value = 0.9842542228653065
$ws.Range("M5").Value = 1.675840666666667
$ws.Range("N5").Value = 5.027522
$ws.Range("O5").Value = 0.03808750486214892
$ws.Range("P5").Value = 0.03808750486214892
$ws.Range("Q5").Value = 981.0723870715206
$ws.Range("R5").Value = 8829.651483643685
$ws.Range("S5").Value = 0.03748778749897296
$ws.Range("T5").Value = 0.03748778749897296

# Row 6
$ws.Range("I6").Value = 0.9842542228653065
$ws.Range("J6").Value = 0.9842542228653065
$ws.Range("O6").Value = 0.096040539564286
$ws.Range("P6").Value = 0.09604053956428601
$ws.Range("S6").Value = 0.09452830663241105
$ws.Range("T6").Value = 0.09452830663241106

# Row 7
$ws.Range("I7").Value = 0.9842542228653065
$ws.Range("J7").Value = 0.9842542228653065
$ws.Range("M7").Value = 38.098149
$ws.Range("N7").Value = 114.294447
$ws.Range("O7").Value = 0.8658719555735651
$ws.Range("P7").Value = 0.865871955573565
$ws.Range("Q7").Value = 22303.45803505373
$ws.Range("R7").Value = 200731.1223154835
$ws.Range("S7").Value = 0.8522381287339226
$ws.Range("T7").Value = 0.8522381287339225

# Row 8
$ws.Range("G8").Value = 9.127409
$ws.Range("H8").Value = 27.382227
$ws.Range("I8").Value = 0.01534569178080505
$ws.Range("J8").Value = 0.01534569178080505
$ws.Range("M8").Value = 1.675840666666667
$ws.Range("N8").Value = 5.027522
$ws.Range("O8").Value = 0.03808750486214892
$ws.Range("P8").Value = 0.03808750486214892
$ws.Range("Q8").Value = 15.29608318349933
$ws.Range("R8").Value = 137.664748651494
$ws.Range("S8").Value = 0.0005844791103144509
$ws.Range("T8").Value = 0.0005844791103144509

# Row 9
$ws.Range("G9").Value = 9.127409
$ws.Range("H9").Value = 27.382227
$ws.Range("I9").Value = 0.01534569178080505
$ws.Range("J9").Value = 0.01534569178080505
$ws.Range("O9").Value = 0.096040539564286
$ws.Range("P9").Value = 0.09604053956428601
$ws.Range("Q9").Value = 38.57023681337034
$ws.Range("R9").Value = 347.132131320333
$ws.Range("S9").Value = 0.001473808518615746
$ws.Range("T9").Value = 0.001473808518615746

# Row 10
$ws.Range("G10").Value = 9.127409
$ws.Range("H10").Value = 27.382227
$ws.Range("I10").Value = 0.01534569178080505
$ws.Range("J10").Value = 0.01534569178080505
$ws.Range("M10").Value = 38.098149
$ws.Range("N10").Value = 114.294447
$ws.Range("O10").Value = 0.8658719555735651
$ws.Range("P10").Value = 0.865871955573565
$ws.Range("Q10").Value = 347.737388065941
$ws.Range("R10").Value = 3129.636492593469
$ws.Range("S10").Value = 0.01328740415187485
$ws.Range("T10").Value = 0.01328740415187485
